# "upload siswa + create read delete anggota"
# Two new rows (15, 16) that previously had no "No" value in column A now get
# sequential numbers (10, 11). Every row below that already had a "No" number
# shifts its number up by 2 to keep the sequence contiguous (A17: 10->12,
# A19:11->13, ... A39:29->31). Finally the active selection moves to A40
# (the first empty row after the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly numbered rows.
$ws.Range("A15").Value = 10
$ws.Range("A16").Value = 11

# Renumber the rest of the "No" column, working from the bottom up so we
# never overwrite a value before it has been read.
$ws.Range("A39").Value = 31
$ws.Range("A37").Value = 30
$ws.Range("A36").Value = 29
$ws.Range("A35").Value = 28
$ws.Range("A34").Value = 27
$ws.Range("A33").Value = 26
$ws.Range("A32").Value = 25
$ws.Range("A31").Value = 24
$ws.Range("A30").Value = 23
$ws.Range("A28").Value = 22
$ws.Range("A27").Value = 21
$ws.Range("A26").Value = 20
$ws.Range("A25").Value = 19
$ws.Range("A24").Value = 18
$ws.Range("A23").Value = 17
$ws.Range("A22").Value = 16
$ws.Range("A21").Value = 15
$ws.Range("A20").Value = 14
$ws.Range("A19").Value = 13
$ws.Range("A17").Value = 12

# Move the active selection to A40 (below the table), matching the saved
# workbook's cursor position after the edit.
$ws.Range("A40").Select()
